$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the
#    H1 title at the top of the document.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Abby and The Witch for Free -
#    Slot Game Review") right before the final "Prompt: ..." paragraph
#    (i.e. right after the last bullet point in "What we don't like").
$lastBullet = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$insertRng = $lastBullet.Range
$insertRng.Collapse(0)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Abby and The Witch for Free - Slot Game Review</w:t></w:r></w:p>'
$insertRng.InsertXML($newParaXml)

# 3. Replace the text of the old "Prompt: ..." paragraph (now the last
#    paragraph) with the meta-description text, keeping its italic run
#    formatting intact.
$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptRng = $promptPara.Range
$oldText = "Prompt: Create a feature image for Abby & The Witch Design a cartoon-style feature image that includes a happy-looking Maya warrior wearing glasses. The image should also incorporate elements from the game " + [char]34 + "Abby & The Witch," + [char]34 + " such as Abby herself, the colorless world, and Baba Yaga's house and cemetery. Use bright colors to contrast the black and white world of the game and make the Maya warrior stand out. Feel free to add other magical elements to the image, like spells, potions, or magical creatures, to give it a more whimsical feel. The image should be eye-catching and convey the spirit of adventure and magic that the game offers to players."
$newText = "Abby and The Witch is an engaging slot game with visually engaging graphics, free spins mode, and low volatility. Read our review and play for free."
$promptRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
